$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 137.5
$ws.Range("I5").Value = 148.9
$ws.Range("J5").Value = 80.5
$ws.Range("K5").Value = 148.9
$ws.Range("L5").Value = 80.5
$ws.Range("M5").Value = -33.90000000000001
$ws.Range("N5").Value = -310.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 94.333336
$ws.Range("I6").Value = 94.333336
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 283.000008
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -171.000008

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 261
$ws.Range("I12").Value = 261
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 261
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -91

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 2420.4546
$ws.Range("I99").Value = 209.71428
$ws.Range("J99").Value = 6289.25
$ws.Range("K99").Value = 629.14284
$ws.Range("L99").Value = 18867.75
$ws.Range("M99").Value = 868.85716
$ws.Range("N99").Value = -21863.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2847.1924
$ws.Range("I100").Value = 2017.3077
$ws.Range("J100").Value = 3677.077
$ws.Range("K100").Value = 2017.3077
$ws.Range("L100").Value = 3677.077
$ws.Range("M100").Value = -1476.3077
$ws.Range("N100").Value = -4759.077

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 597.0909
$ws.Range("I101").Value = 471.875
$ws.Range("J101").Value = 931
$ws.Range("K101").Value = 1415.625
$ws.Range("L101").Value = 2793
$ws.Range("M101").Value = 206.375
$ws.Range("N101").Value = -6037

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 14173.4
$ws.Range("I106").Value = 10216.75
$ws.Range("J106").Value = 30000
$ws.Range("K106").Value = 10216.75
$ws.Range("L106").Value = 30000
$ws.Range("M106").Value = -9585.75
$ws.Range("N106").Value = -31262

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 14773.6
$ws.Range("I113").Value = 19990
$ws.Range("J113").Value = 6949
$ws.Range("K113").Value = 19990
$ws.Range("L113").Value = 6949
$ws.Range("M113").Value = -16736
$ws.Range("N113").Value = -13457

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 2000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 6000
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -16000

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3599.0625
$ws.Range("I132").Value = 3667.4194
$ws.Range("J132").Value = 1480
$ws.Range("K132").Value = 11002.2582
$ws.Range("L132").Value = 4440
$ws.Range("M132").Value = -8472.2582
$ws.Range("N132").Value = -9500

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 80875.17999999999
$ws.Range("I140").Value = 50709
$ws.Range("J140").Value = 98113
$ws.Range("K140").Value = 50709
$ws.Range("L140").Value = 98113
$ws.Range("M140").Value = -45529
$ws.Range("N140").Value = -108473

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21299.932
$ws.Range("I32").Value = 23374.611
$ws.Range("J32").Value = 10004.444
$ws.Range("K32").Value = 23374.611
$ws.Range("L32").Value = 10004.444
$ws.Range("M32").Value = -23087.611
$ws.Range("N32").Value = -10578.444

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3818.2334
$ws.Range("I45").Value = 3051
$ws.Range("J45").Value = 5352.7
$ws.Range("K45").Value = 3051
$ws.Range("L45").Value = 5352.7
$ws.Range("M45").Value = -2674
$ws.Range("N45").Value = -6106.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1954.5
$ws.Range("I61").Value = 1945.5
$ws.Range("J61").Value = 1999.5
$ws.Range("K61").Value = 1945.5
$ws.Range("L61").Value = 1999.5
$ws.Range("M61").Value = -1733.5
$ws.Range("N61").Value = -2423.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3000
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2314

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3000
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 15000
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -11568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 47981.137
$ws.Range("I132").Value = 55076.21
$ws.Range("J132").Value = 3045.6667
$ws.Range("K132").Value = 165228.63
$ws.Range("L132").Value = 9137.000100000001
$ws.Range("M132").Value = -162698.63
$ws.Range("N132").Value = -14197.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1954.5
$ws.Range("I136").Value = 1945.5
$ws.Range("J136").Value = 1999.5
$ws.Range("K136").Value = 5836.5
$ws.Range("L136").Value = 5998.5
$ws.Range("M136").Value = -3286.5
$ws.Range("N136").Value = -11098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 6173.5
$ws.Range("I5").Value = 4231.3335
$ws.Range("J5").Value = 12000
$ws.Range("K5").Value = 4231.3335
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = -4118.3335
$ws.Range("N5").Value = -12226

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 3000
$ws.Range("I7").Value = 3000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2887
$ws.Range("N7").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4339.6
$ws.Range("I20").Value = 4249.5
$ws.Range("J20").Value = 4399.6665
$ws.Range("K20").Value = 4249.5
$ws.Range("L20").Value = 4399.6665
$ws.Range("M20").Value = -4002.5
$ws.Range("N20").Value = -4893.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 1956.25
$ws.Range("I128").Value = 1956.25
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 5868.75
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -3378.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2081.4
$ws.Range("I134").Value = 1539
$ws.Range("J134").Value = 5607
$ws.Range("K134").Value = 4617
$ws.Range("L134").Value = 16821
$ws.Range("M134").Value = -2082
$ws.Range("N134").Value = -21891

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 767.375
$ws.Range("I16").Value = 939.8333
$ws.Range("J16").Value = 250
$ws.Range("K16").Value = 939.8333
$ws.Range("L16").Value = 250
$ws.Range("M16").Value = -652.8333
$ws.Range("N16").Value = -824

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1608.2
$ws.Range("I22").Value = 247
$ws.Range("J22").Value = 3163.8572
$ws.Range("K22").Value = 247
$ws.Range("L22").Value = 3163.8572
$ws.Range("M22").Value = 103
$ws.Range("N22").Value = -3863.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3641.4
$ws.Range("I99").Value = 4121
$ws.Range("J99").Value = 3435.8572
$ws.Range("K99").Value = 4121
$ws.Range("L99").Value = 3435.8572
$ws.Range("M99").Value = -2623
$ws.Range("N99").Value = -6431.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 767.375
$ws.Range("I113").Value = 939.8333
$ws.Range("J113").Value = 250
$ws.Range("K113").Value = 939.8333
$ws.Range("L113").Value = 250
$ws.Range("M113").Value = 1230.1667
$ws.Range("N113").Value = -4590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2748.5
$ws.Range("I122").Value = 2738.7
$ws.Range("J122").Value = 2797.5
$ws.Range("K122").Value = 8216.099999999999
$ws.Range("L122").Value = 8392.5
$ws.Range("M122").Value = -5766.099999999999
$ws.Range("N122").Value = -13292.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3641.4
$ws.Range("I126").Value = 4121
$ws.Range("J126").Value = 3435.8572
$ws.Range("K126").Value = 12363
$ws.Range("L126").Value = 10307.5716
$ws.Range("M126").Value = -9893
$ws.Range("N126").Value = -15247.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2010448.9
$ws.Range("I4").Value = 1387811.1
$ws.Range("J4").Value = 4501000
$ws.Range("K4").Value = 4163433.3
$ws.Range("L4").Value = 13503000
$ws.Range("M4").Value = -4163321.3
$ws.Range("N4").Value = -13503224

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1222.5555
$ws.Range("I5").Value = 1156.625
$ws.Range("J5").Value = 1750
$ws.Range("K5").Value = 3469.875
$ws.Range("L5").Value = 5250
$ws.Range("M5").Value = -3357.875
$ws.Range("N5").Value = -5474

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 223
$ws.Range("I12").Value = 223.16667
$ws.Range("J12").Value = 222.90909
$ws.Range("K12").Value = 669.50001
$ws.Range("L12").Value = 668.72727
$ws.Range("M12").Value = -496.50001
$ws.Range("N12").Value = -1014.72727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 7126.647
$ws.Range("I55").Value = 4244
$ws.Range("J55").Value = 8013.615
$ws.Range("K55").Value = 12732
$ws.Range("L55").Value = 24040.845
$ws.Range("M55").Value = -12555
$ws.Range("N55").Value = -24394.845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7665.5
$ws.Range("I56").Value = 7665.5
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 7665.5
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -7135.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 4880
$ws.Range("I134").Value = 4880
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 14640
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -9570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1222.5555
$ws.Range("I135").Value = 1156.625
$ws.Range("J135").Value = 1750
$ws.Range("K135").Value = 10409.625
$ws.Range("L135").Value = 15750
$ws.Range("M135").Value = -7874.625
$ws.Range("N135").Value = -20820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3332.6667
$ws.Range("I102").Value = 1999.5
$ws.Range("J102").Value = 5999
$ws.Range("K102").Value = 1999.5
$ws.Range("L102").Value = 5999
$ws.Range("M102").Value = -377.5
$ws.Range("N102").Value = -9243

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3473.182
$ws.Range("I113").Value = 2974.5715
$ws.Range("J113").Value = 4345.75
$ws.Range("K113").Value = 2974.5715
$ws.Range("L113").Value = 4345.75
$ws.Range("M113").Value = -804.5715
$ws.Range("N113").Value = -8685.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6406
$ws.Range("I126").Value = 3846.75
$ws.Range("J126").Value = 7685.625
$ws.Range("K126").Value = 11540.25
$ws.Range("L126").Value = 23056.875
$ws.Range("M126").Value = -9070.25
$ws.Range("N126").Value = -27996.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 48577.637
$ws.Range("I132").Value = 61900.53
$ws.Range("J132").Value = 3279.8
$ws.Range("K132").Value = 185701.59
$ws.Range("L132").Value = 9839.400000000001
$ws.Range("M132").Value = -183171.59
$ws.Range("N132").Value = -14899.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3428.9412
$ws.Range("I61").Value = 2986.2
$ws.Range("J61").Value = 6749.5
$ws.Range("K61").Value = 2986.2
$ws.Range("L61").Value = 6749.5
$ws.Range("M61").Value = -2784.2
$ws.Range("N61").Value = -7153.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2208.9143
$ws.Range("I93").Value = 2095.5
$ws.Range("J93").Value = 2284.524
$ws.Range("K93").Value = 2095.5
$ws.Range("L93").Value = 2284.524
$ws.Range("M93").Value = -847.5
$ws.Range("N93").Value = -4780.523999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3428.9412
$ws.Range("I113").Value = 2986.2
$ws.Range("J113").Value = 6749.5
$ws.Range("K113").Value = 2986.2
$ws.Range("L113").Value = 6749.5
$ws.Range("M113").Value = -816.1999999999998
$ws.Range("N113").Value = -11089.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 28256.596
$ws.Range("I132").Value = 32226.4
$ws.Range("J132").Value = 5572
$ws.Range("K132").Value = 96679.20000000001
$ws.Range("L132").Value = 16716
$ws.Range("M132").Value = -94149.20000000001
$ws.Range("N132").Value = -21776

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6336.3335
$ws.Range("I136").Value = 5603.6
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 16810.8
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -14260.8
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 3099.75
$ws.Range("I5").Value = 3999
$ws.Range("J5").Value = 2800
$ws.Range("K5").Value = 3999
$ws.Range("L5").Value = 2800
$ws.Range("M5").Value = -3887
$ws.Range("N5").Value = -3024

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2064.1428
$ws.Range("I81").Value = 1816.3334
$ws.Range("J81").Value = 2250
$ws.Range("K81").Value = 3632.6668
$ws.Range("L81").Value = 4500
$ws.Range("M81").Value = -2571.6668
$ws.Range("N81").Value = -6622

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2064.1428
$ws.Range("I84").Value = 1816.3334
$ws.Range("J84").Value = 2250
$ws.Range("K84").Value = 18163.334
$ws.Range("L84").Value = 22500
$ws.Range("M84").Value = -12859.334
$ws.Range("N84").Value = -33108

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 737.0909
$ws.Range("I100").Value = 646.1111
$ws.Range("J100").Value = 1146.5
$ws.Range("K100").Value = 1292.2222
$ws.Range("L100").Value = 2293
$ws.Range("M100").Value = -751.2221999999999
$ws.Range("N100").Value = -3375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 66438.164
$ws.Range("I132").Value = 70862.32000000001
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 212586.96
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -210056.96
$ws.Range("N132").Value = -18560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1813.3793
$ws.Range("I136").Value = 1838.4762
$ws.Range("J136").Value = 1747.5
$ws.Range("K136").Value = 5515.4286
$ws.Range("L136").Value = 5242.5
$ws.Range("M136").Value = -2965.4286
$ws.Range("N136").Value = -10342.5
